$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.301.58'
$ws.Range("E2").Value = '  +1.45%  '
$ws.Range("D3").Value = '1.892.89'
$ws.Range("E3").Value = '  +1.38%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.49'
$ws.Range("E5").Value = '  +1.24%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.690'
$ws.Range("E6").Value = '  +3.05%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '42.77'
$ws.Range("E8").Value = '  +1.68%  '
$ws.Range("E9").Value = '  +5.45%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '56.22'
$ws.Range("E10").Value = '  +10.69%  '
$ws.Range("E11").Value = '  +2.13%  '
$ws.Range("E12").Value = '  +1.50%  '
$ws.Range("E13").Value = '  +9.05%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.796'
$ws.Range("E14").Value = '  +12.61%  '
$ws.Range("D15").Value = '2.171.19'
$ws.Range("E15").Value = '  +1.58%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.00'
$ws.Range("E16").Value = '  +4.14%  '
$ws.Range("D17").Value = '1.904.19'
$ws.Range("E17").Value = '  +2.10%  '
$ws.Range("D18").Value = '35.384.78'
$ws.Range("E18").Value = '  +1.78%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '73.44'
$ws.Range("E19").Value = '  +2.07%  '
$ws.Range("D20").Value = '0.0₃0827'
$ws.Range("E20").Value = '  +2.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '243.77'
$ws.Range("E21").Value = '  +0.81%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.91'
$ws.Range("E22").Value = '  +3.23%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.21'
$ws.Range("E23").Value = '  +7.25%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.69'
$ws.Range("E24").Value = '  +8.64%  '
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.17'
$ws.Range("E26").Value = '  +1.72%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '166.70'
$ws.Range("E27").Value = '  +2.37%  '
$ws.Range("E28").Value = '  +3.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.31'
$ws.Range("E29").Value = '  +1.88%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.127'
$ws.Range("E30").Value = '  +2.14%  '
$ws.Range("E31").Value = '  +6.72%  '
$ws.Range("E32").Value = '  +4.94%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.23'
$ws.Range("E33").Value = '  +3.20%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.85'
$ws.Range("E34").Value = '  +24.23%  '
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("E36").Value = '  -14.91%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.850'
$ws.Range("E37").Value = '  +2.96%  '
$ws.Range("E38").Value = '  +1.91%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0720'
$ws.Range("E39").Value = '  +8.58%  '
$ws.Range("E40").Value = '  +7.72%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '98.83'
$ws.Range("E41").Value = '  +1.84%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '16.98'
$ws.Range("E42").Value = '  +0.77%  '
$ws.Range("E43").Value = '  +1.37%  '
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '1.336.47'
$ws.Range("E44").Value = '  +4.61%  '
$ws.Range("B45").Value = 'Gas'
$ws.Range("C45").Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.52'
$ws.Range("E45").Value = '  +14.03%  '
$ws.Range("E46").Value = '  +3.38%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0809'
$ws.Range("E47").Value = '  -1.91%  '
$ws.Range("E48").Value = '  +0.73%  '
$ws.Range("E49").Value = '  +0.43%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.35'
$ws.Range("E50").Value = '  +1.84%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '42.53'
$ws.Range("E51").Value = '  +0.58%  '
